$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Relabel the "D0" row header to "D0 [cm/s]"
$ws1.Range("C8").Value = "D0 [cm/s]"

# Update the D0 value used in the D-column diffusivity formula (was 3e10, now 3e7)
$ws1.Range("D8").Value = 30000000

# Update the active cell selection on Sheet1
$ws1.Range("F16").Select()
